# Auto-generated Excel COM-interop script
# Applies the cryptocurrency price/volume updates described by the commit diff
# to the active worksheet, preserving each cell's original text (inline string) type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, and non-numeric-looking price/volume strings).
# These are safe to assign directly since Excel will not reinterpret them as numbers.
$plainUpdates = @{
    'D2' = '42.250.34'
    'E2' = '  -1.57%  '
    'D3' = '2.517.41'
    'E3' = '  -2.02%  '
    'E4' = '  -0.17%  '
    'E5' = '  +0.18%  '
    'E6' = '  -0.73%  '
    'E7' = '  +1.66%  '
    'E8' = '  +0.04%  '
    'E9' = '  -2.66%  '
    'E10' = '  -0.55%  '
    'E11' = '  -0.31%  '
    'E12' = '  -0.10%  '
    'E13' = '  -1.76%  '
    'D14' = '2.897.33'
    'E14' = '  -2.19%  '
    'D15' = '2.548.74'
    'E15' = '  +0.50%  '
    'E16' = '  +4.75%  '
    'E17' = '  -3.01%  '
    'D18' = '42.300.42'
    'E18' = '  -1.40%  '
    'E19' = '  +0.13%  '
    'D20' = '0.0₃0969'
    'E20' = '  -2.70%  '
    'E21' = '  -3.02%  '
    'E22' = '  -1.51%  '
    'E23' = '  -1.68%  '
    'E24' = '  -1.99%  '
    'E25' = '  -5.30%  '
    'E26' = '  -6.63%  '
    'E27' = '  +0.22%  '
    'E28' = '  +10.01%  '
    'E29' = '  -0.30%  '
    'E30' = '  +1.11%  '
    'E31' = '  -1.79%  '
    'E32' = '  -0.44%  '
    'E33' = '  -2.20%  '
    'E34' = '  -2.75%  '
    'B35' = 'ARBITRUM'
    'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'E35' = '  -4.75%  '
    'B36' = 'Celestia'
    'C36' = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
    'E36' = '  +1.38%  '
    'B37' = 'WEMIXToken'
    'C37' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'E37' = '  -5.17%  '
    'E38' = '  +0.32%  '
    'E39' = '  +3.32%  '
    'E40' = '  -0.98%  '
    'E41' = '  -1.72%  '
    'E42' = '  -1.18%  '
    'B43' = 'ApeXProtocol'
    'C43' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'E43' = '  -1.46%  '
    'B44' = 'FirstDigitalUSD'
    'C44' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'E44' = '  -0.25%  '
    'B45' = 'VeChain'
    'C45' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E45' = '  -3.69%  '
    'B46' = 'Maker'
    'C46' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D46' = '2.031.58'
    'E46' = '  -1.95%  '
    'E47' = '  -1.16%  '
    'E48' = '  -3.27%  '
    'D49' = '2.758.19'
    'E49' = '  -2.21%  '
    'E50' = '  -4.64%  '
    'E51' = '  -1.56%  '
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Updates whose new text would otherwise be auto-converted to a number by Excel
# (e.g. '303.12' or '0.997'). Force the cell to Text format first so the original
# string representation (with its exact decimal formatting) is preserved verbatim.
$textForcedUpdates = @{
    'D4' = '0.997'
    'D5' = '303.12'
    'D6' = '96.28'
    'D7' = '0.585'
    'D10' = '36.45'
    'D11' = '0.0806'
    'D12' = '7.63'
    'D16' = '15.03'
    'D17' = '0.860'
    'D19' = '12.90'
    'D21' = '6.43'
    'D22' = '70.90'
    'D23' = '249.97'
    'D24' = '2.90'
    'D25' = '2.01'
    'D26' = '26.91'
    'D28' = '2.33'
    'D29' = '10.25'
    'D30' = '37.87'
    'D31' = '5.92'
    'D32' = '154.66'
    'D33' = '3.31'
    'D34' = '0.0785'
    'D35' = '2.06'
    'D36' = '18.60'
    'D37' = '2.61'
    'D39' = '24.18'
    'D41' = '3.37'
    'D42' = '3.82'
    'D43' = '2.02'
    'D44' = '0.996'
    'D45' = '0.0299'
    'D47' = '84.42'
    'D50' = '101.55'
}

foreach ($cellRef in $textForcedUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$cellRef]
}
